# Apply updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "94.742.82"
$ws.Range("E2").Value = "  +2.53%  "
$ws.Range("D3").Value = "3.120.33"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "614.52"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +2.20%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.393"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D11").Value = "3.118.87"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  -1.92%  "
$ws.Range("E13").Value = "  -2.41%  "
$ws.Range("D14").Value = "94.350.90"
$ws.Range("E14").Value = "  +2.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.77"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.29"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -4.07%  "
$ws.Range("D17").Value = "3.696.62"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "3.120.50"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.68"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.01"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.69%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.65%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "450.97"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000199"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.01"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.33"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +5.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.64"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.13"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.96%  "
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "85.98"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "3.287.90"
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.258"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +9.26%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.183"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +9.26%  "
$ws.Range("E33").Value = "  -10.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.38"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.989"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.86%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.95"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("E37").Value = "  -0.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.02"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.459"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +5.57%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.30"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "477.64"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.60%  "
$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.98"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.62"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -14.13%  "
$ws.Range("E45").Value = "  -7.60%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "161.59"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.695"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.89"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.44%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0325"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.81%  "
